$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new "through" date
$wb.Worksheets.Item(1).Name = "Through 2021-12-23"

# Update the row label for December to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-23)"

# Update December (row 13) values for each year column (B..H = 2015..2021)
$ws.Range("B13").Value = 32
$ws.Range("C13").Value = 73
$ws.Range("D13").Value = 92
$ws.Range("E13").Value = 53
$ws.Range("F13").Value = 47
$ws.Range("G13").Value = 112
$ws.Range("H13").Value = 152

# Update Total (row 14) values for each year column (B..H = 2015..2021)
$ws.Range("B14").Value = 323
$ws.Range("C14").Value = 636
$ws.Range("D14").Value = 913
$ws.Range("E14").Value = 735
$ws.Range("F14").Value = 581
$ws.Range("G14").Value = 1376
$ws.Range("H14").Value = 1795
